{"js": "// 1) \"... analyze the data using a Stacked Column Pivot Chart based on\n//    sub-category, play...\" -> \"... analyzed the data using a Stacked\n//    Column Pivot Chart based on sub-category, play...\"\nconst analyzeResults = context.document.body.search(\n  \"analyze the data using a Stacked Column Pivot Chart based on sub-category\",\n  { matchCase: true }\n);\nanalyzeResults.load(\"text\");\nawait context.sync();\n\nif (analyzeResults.items.length > 0) {\n  analyzeResults.items[0].insertText(\n    \"analyzed the data using a Stacked Column Pivot Chart based on sub-category\",\n    \"Replace\"\n  );\n}\n\n// 2) \"Some limitations of this data set?\" -> \"Some limitations of data set?\"\nconst limitationsResults = context.document.body.search(\n  \"limitations of this data\",\n  { matchCase: true }\n);\nlimitationsResults.load(\"text\");\nawait context.sync();\n\nif (limitationsResults.items.length > 0) {\n  limitationsResults.items[0].insertText(\"limitations of data\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"analyze the data\" -> \"analyzed the data\"\n#    (the sub-category sentence; the category sentence already reads \"analyzed\")\n$rng1 = $d.Content\n$rng1.Find.Execute(\n    \"analyze the data using a Stacked Column Pivot Chart based on sub-category\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"analyzed the data using a Stacked Column Pivot Chart based on sub-category\",\n    2\n)\n\n# 2) \"ome limitations of this data\" -> \"ome limitations of data\"\n$rng2 = $d.Content\n$rng2.Find.Execute(\n    \"limitations of this data\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"limitations of data\",\n    2\n)\n"}
